$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H98").Value = 1757.2307
$ws.Range("I98").Value = 1702.3334
$ws.Range("K98").Value = 1702.3334
$ws.Range("M98").Value = -204.3334
$ws.Range("H111").Value = 962.1667
$ws.Range("J111").Value = 1566.4
$ws.Range("L111").Value = 4699.200000000001
$ws.Range("N111").Value = -10833.2
$ws.Range("H122").Value = 1757.2307
$ws.Range("I122").Value = 1702.3334
$ws.Range("K122").Value = 5107.0002
$ws.Range("M122").Value = -2657.0002
$ws.Range("H125").Value = 5535.8335
$ws.Range("I125").Value = 4220.909
$ws.Range("K125").Value = 37988.181
$ws.Range("M125").Value = -35528.181
$ws.Range("H138").Value = 3064.3062
$ws.Range("I138").Value = 2881.6072
$ws.Range("K138").Value = 8644.821599999999
$ws.Range("M138").Value = -3504.821599999999

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 2491.9412
$ws.Range("I2").Value = 885.6667
$ws.Range("J2").Value = 8687.571
$ws.Range("K2").Value = 885.6667
$ws.Range("L2").Value = 8687.571
$ws.Range("M2").Value = -772.6667
$ws.Range("N2").Value = -8913.571
$ws.Range("H32").Value = 14072.622
$ws.Range("I32").Value = 5330.093
$ws.Range("K32").Value = 5330.093
$ws.Range("M32").Value = -5043.093
$ws.Range("H61").Value = 13111.394
$ws.Range("I61").Value = 7196.8
$ws.Range("J61").Value = 22210.77
$ws.Range("K61").Value = 7196.8
$ws.Range("L61").Value = 22210.77
$ws.Range("M61").Value = -6984.8
$ws.Range("N61").Value = -22634.77
$ws.Range("H88").Value = 2199.8572
$ws.Range("I88").Value = 1624.75
$ws.Range("K88").Value = 1624.75
$ws.Range("M88").Value = -1218.75
$ws.Range("H91").Value = 2199.8572
$ws.Range("I91").Value = 1624.75
$ws.Range("K91").Value = 1624.75
$ws.Range("M91").Value = -220.75
$ws.Range("H110").Value = 7497.6
$ws.Range("I110").Value = 4888.3335
$ws.Range("J110").Value = 17934.666
$ws.Range("K110").Value = 4888.3335
$ws.Range("L110").Value = 17934.666
$ws.Range("M110").Value = -2843.3335
$ws.Range("N110").Value = -22024.666
$ws.Range("H116").Value = 2491.9412
$ws.Range("I116").Value = 885.6667
$ws.Range("J116").Value = 8687.571
$ws.Range("K116").Value = 885.6667
$ws.Range("L116").Value = 8687.571
$ws.Range("M116").Value = 1408.3333
$ws.Range("N116").Value = -13275.571
$ws.Range("H122").Value = 3265.25
$ws.Range("I122").Value = 2410.6875
$ws.Range("K122").Value = 7232.0625
$ws.Range("M122").Value = -4782.0625
$ws.Range("H136").Value = 13111.394
$ws.Range("I136").Value = 7196.8
$ws.Range("J136").Value = 22210.77
$ws.Range("K136").Value = 21590.4
$ws.Range("L136").Value = 66632.31
$ws.Range("M136").Value = -19040.4
$ws.Range("N136").Value = -71732.31
$ws.Range("H139").Value = 73644.22
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 2491.9412
$ws.Range("I3").Value = 885.6667
$ws.Range("J3").Value = 8687.571
$ws.Range("K3").Value = 885.6667
$ws.Range("L3").Value = 8687.571
$ws.Range("M3").Value = -771.6667
$ws.Range("N3").Value = -8915.571
$ws.Range("H80").Value = 1066
$ws.Range("I80").Value = 281.57144
$ws.Range("J80").Value = 1304.7391
$ws.Range("K80").Value = 281.57144
$ws.Range("L80").Value = 1304.7391
$ws.Range("M80").Value = 716.4285600000001
$ws.Range("N80").Value = -3300.7391
$ws.Range("H83").Value = 1066
$ws.Range("I83").Value = 281.57144
$ws.Range("J83").Value = 1304.7391
$ws.Range("K83").Value = 1407.8572
$ws.Range("L83").Value = 6523.6955
$ws.Range("M83").Value = 3584.1428
$ws.Range("N83").Value = -16507.6955
$ws.Range("H86").Value = 3870.8333
$ws.Range("I86").Value = 3642.2222
$ws.Range("J86").Value = 4099.4443
$ws.Range("K86").Value = 3642.2222
$ws.Range("L86").Value = 4099.4443
$ws.Range("M86").Value = -2519.2222
$ws.Range("N86").Value = -6345.4443
$ws.Range("H89").Value = 3870.8333
$ws.Range("I89").Value = 3642.2222
$ws.Range("J89").Value = 4099.4443
$ws.Range("K89").Value = 18211.111
$ws.Range("L89").Value = 20497.2215
$ws.Range("M89").Value = -12595.111
$ws.Range("N89").Value = -31729.2215

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H122").Value = 8395.25
$ws.Range("I122").Value = 6805.625
$ws.Range("J122").Value = 11574.5
$ws.Range("K122").Value = 20416.875
$ws.Range("L122").Value = 34723.5
$ws.Range("M122").Value = -17966.875
$ws.Range("N122").Value = -39623.5
$ws.Range("H132").Value = 5841.5264
$ws.Range("I132").Value = 2120.7727
$ws.Range("K132").Value = 6362.3181
$ws.Range("M132").Value = -3832.3181

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H14").Value = 215.75
$ws.Range("I14").Value = 215.75
$ws.Range("K14").Value = 647.25
$ws.Range("M14").Value = -474.25
$ws.Range("H117").Value = 1896.2106
$ws.Range("I117").Value = 450
$ws.Range("J117").Value = 2412.7144
$ws.Range("K117").Value = 1350
$ws.Range("L117").Value = 7238.1432
$ws.Range("M117").Value = 2092
$ws.Range("N117").Value = -14122.1432
$ws.Range("H129").Value = 1619
$ws.Range("I129").Value = 1146.8334
$ws.Range("J129").Value = 2563.3333
$ws.Range("K129").Value = 3440.5002
$ws.Range("L129").Value = 7689.999899999999
$ws.Range("M129").Value = 1559.4998
$ws.Range("N129").Value = -17689.9999
$ws.Range("H132").Value = 1277.4286
$ws.Range("I132").Value = 1466.1111
$ws.Range("J132").Value = 937.8
$ws.Range("K132").Value = 13194.9999
$ws.Range("L132").Value = 8440.199999999999
$ws.Range("M132").Value = -10664.9999
$ws.Range("N132").Value = -13500.2

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H102").Value = 4385.4736
$ws.Range("I102").Value = 2119.6
$ws.Range("J102").Value = 6903.1113
$ws.Range("K102").Value = 2119.6
$ws.Range("L102").Value = 6903.1113
$ws.Range("M102").Value = -497.5999999999999
$ws.Range("N102").Value = -10147.1113
$ws.Range("H122").Value = 5201.269
$ws.Range("I122").Value = 2826.8
$ws.Range("K122").Value = 8480.400000000001
$ws.Range("M122").Value = -6030.400000000001
$ws.Range("H132").Value = 8795.049999999999
$ws.Range("I132").Value = 3825.1538
$ws.Range("J132").Value = 18024.857
$ws.Range("K132").Value = 11475.4614
$ws.Range("L132").Value = 54074.571
$ws.Range("M132").Value = -8945.4614
$ws.Range("N132").Value = -59134.571

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 7706.1763
$ws.Range("I7").Value = 6002.3335
$ws.Range("K7").Value = 6002.3335
$ws.Range("M7").Value = -5890.3335
$ws.Range("H40").Value = 10662.895
$ws.Range("I40").Value = 7100.6665
$ws.Range("K40").Value = 7100.6665
$ws.Range("M40").Value = -6964.6665
$ws.Range("H68").Value = 49616.332
$ws.Range("J68").Value = 49499.5
$ws.Range("L68").Value = 49499.5
$ws.Range("N68").Value = -50997.5
$ws.Range("H71").Value = 49616.332
$ws.Range("J71").Value = 49499.5
$ws.Range("L71").Value = 247497.5
$ws.Range("N71").Value = -254985.5
$ws.Range("H122").Value = 9784.315000000001
$ws.Range("I122").Value = 7555.4443
$ws.Range("K122").Value = 22666.3329
$ws.Range("M122").Value = -20216.3329
$ws.Range("H126").Value = 7706.1763
$ws.Range("I126").Value = 6002.3335
$ws.Range("K126").Value = 18007.0005
$ws.Range("M126").Value = -15537.0005

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H62").Value = 3099.5715
$ws.Range("I62").Value = 2829.4
$ws.Range("K62").Value = 2829.4
$ws.Range("M62").Value = -2205.4
$ws.Range("H65").Value = 3099.5715
$ws.Range("I65").Value = 2829.4
$ws.Range("K65").Value = 14147
$ws.Range("M65").Value = -11027
$ws.Range("H100").Value = 527.0952
$ws.Range("I100").Value = 450.66666
$ws.Range("K100").Value = 901.33332
$ws.Range("M100").Value = -360.33332
$ws.Range("H126").Value = 22558.143
$ws.Range("I126").Value = 22090.63
$ws.Range("J126").Value = 26999.5
$ws.Range("K126").Value = 66271.89
$ws.Range("L126").Value = 80998.5
$ws.Range("M126").Value = -63801.89
$ws.Range("N126").Value = -85938.5
